$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("D3").Value = 11
$ws.Range("D4").Value = 11
$ws.Range("D5").Value = 12
$ws.Range("D6").Value = 13
$ws.Range("D8").Value = 14
$ws.Range("D9").Value = 15
$ws.Range("D10").Value = 16
$ws.Range("D11").Value = 14
$ws.Range("D12").Value = 17
$ws.Range("D13").Value = 18
$ws.Range("D14").Value = 19
$ws.Range("D15").Value = 20
$ws.Range("D16").Value = 21
$ws.Range("D17").Value = 22
$ws.Range("D18").Value = 25
$ws.Range("D19").Value = 24
$ws.Range("D20").Value = 23
$ws.Range("D21").Value = 26
$ws.Range("D22").Value = 28
$ws.Range("D23").Value = 29
$ws.Range("D24").Value = 30
$ws.Range("D25").Value = 30
$ws.Range("D26").Value = 31
$ws.Range("D27").Value = 30
$ws.Range("D28").Value = 32
$ws.Range("D29").Value = 33
$ws.Range("D30").Value = 31
$ws.Range("D31").Value = 32
$ws.Range("D32").Value = 30
$ws.Range("D33").Value = 34
$ws.Range("D34").Value = 35
$ws.Range("D35").Value = 36
$ws.Range("D36").Value = 38
$ws.Range("D37").Value = 38
$ws.Range("D38").Value = 38
$ws.Range("D39").Value = 39
$ws.Range("D40").Value = 39
$ws.Range("D41").Value = 39
$ws.Range("D42").Value = 40
$ws.Range("D43").Value = 41
$ws.Range("D44").Value = 42
$ws.Range("D45").Value = 41
$ws.Range("D46").Value = 43
$ws.Range("D47").Value = 45
$ws.Range("D48").Value = 44
$ws.Range("D49").Value = 44
$ws.Range("D50").Value = 44
$ws.Range("D51").Value = 46
$ws.Range("D52").Value = 48

$ws.Range("D53").Select()
